# Weekly update: insert the latest week's "Cebollín" price record right
# after the existing header/first-week row (new row 3), pushing every
# previously recorded week down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank row at row 3; everything below (old rows 3..38)
# shifts down to become rows 4..39.
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with this week's record.
$ws.Cells.Item(3, 1).Value  = 11
$ws.Cells.Item(3, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(3, 3).Value  = "Bíobío"
$ws.Cells.Item(3, 4).Value  = 44687
$ws.Cells.Item(3, 5).Value  = 8
$ws.Cells.Item(3, 6).Value  = 100112037
$ws.Cells.Item(3, 7).Value  = "Cebollín"
$ws.Cells.Item(3, 8).Value  = "Sin especificar"
$ws.Cells.Item(3, 9).Value  = "Primera"
$ws.Cells.Item(3, 10).Value = 220
$ws.Cells.Item(3, 11).Value = 8000
$ws.Cells.Item(3, 12).Value = 8500
$ws.Cells.Item(3, 13).Value = 8273
$ws.Cells.Item(3, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(3, 15).Value = "Región Metropolitana"
$ws.Cells.Item(3, 16).Value = 230
$ws.Cells.Item(3, 17).Value = 36
$ws.Cells.Item(3, 18).Value = "Hortaliza"
